# Update bug report OLX: rebuild sheet1 ("Лист1") as an English test-case
# report template (previously a Russian OLX bug report).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Reset the existing layout: drop old merges, wipe all cell content
#    and formatting so the sheet starts from a clean slate (this also
#    drops now-unused shared strings on save).
# ---------------------------------------------------------------------
$ws.Range("A1:A6").UnMerge()
$ws.Range("D5:D6").UnMerge()
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# 2. Write the new cell text, in the same order the cells are first
#    populated in the target workbook, so the shared-string table is
#    rebuilt in that same order.
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 'URL: https://www.olx.ua'
$ws.Range("B14").Value = '1. Кликнуть на "Мой профиль".'
$ws.Range("B15").Value = '2. Для изменеия языка на сайте вверху на ленте кликнуть на "мова". '
$ws.Range("D3").Value = 'Test Designed date:'
$ws.Range("D4").Value = 'Test Executed by:'
$ws.Range("D5").Value = 'Test Executed date:'
$ws.Range("A1").Value = 'Module name: Registration of a new customer'
$ws.Range("A3").Value = 'Test Priority:'
$ws.Range("A4").Value = 'Tesk?AC ID: TC1'
$ws.Range("A5").Value = 'Test Title: Registration of a new customer with valid data'
$ws.Range("A6").Value = 'Description:'
$ws.Range("A8").Value = 'Pre-conditions:'
$ws.Range("A9").Value = 'Data for testing:  Vasya Pupkin@gmail.com'
$ws.Range("B13").Value = 'Step'
$ws.Range("D6").Value = 'Enwiroment: Windows 10, Chrome'
$ws.Range("D2").Value = 'Test Designed by: Lutsenko D.'
$ws.Range("C13").Value = 'Expected Result'
$ws.Range("D13").Value = 'Actual Result'
$ws.Range("C1").Value = 'The word "Текущий" is not translated correctly. The word "streaming" with the letter "ы" is displayed.'
$ws.Range("C14").Value = 'The registration page opens. The entire page is displayed translated into Ukrainian.'
$ws.Range("D14").Value = 'The word "Текущий" is not translated correctly. The word "streaming" with the letter "ы" is displayed.'

$ws.Range("A14").Value = 1
$ws.Range("A15").Value = 2

# ---------------------------------------------------------------------
# 3. Alignment / formatting per cell (mirrors the cellXfs used in the
#    target file).
# ---------------------------------------------------------------------
function Set-CellAlign($rng, $h, $v, $wrap) {
    if ($h -ne $null) { $rng.HorizontalAlignment = $h }
    if ($v -ne $null) { $rng.VerticalAlignment = $v }
    $rng.WrapText = $wrap
}

$xlLeft = -4131
$xlCenter = -4108
$xlRight = -4152
$xlTop = -4160
$xlVCenter = -4108
$xlGeneralH = 1
$xlGeneralV = -4107

Set-CellAlign $ws.Range("A14") $xlCenter $xlVCenter $false
Set-CellAlign $ws.Range("B14") $xlLeft $xlTop $false
Set-CellAlign $ws.Range("C14") $xlGeneralH $xlVCenter $true
Set-CellAlign $ws.Range("D14") $xlLeft $xlVCenter $true

Set-CellAlign $ws.Range("A15") $xlCenter $xlVCenter $true
Set-CellAlign $ws.Range("B15") $xlLeft $xlTop $true
Set-CellAlign $ws.Range("C15") $xlGeneralH $xlVCenter $true
Set-CellAlign $ws.Range("D15") $xlLeft $xlVCenter $true

Set-CellAlign $ws.Range("A16") $xlCenter $xlVCenter $true
Set-CellAlign $ws.Range("B16") $xlRight $xlVCenter $false

Set-CellAlign $ws.Range("A17") $xlCenter $xlVCenter $true

Set-CellAlign $ws.Range("A18") $xlCenter $xlVCenter $true
Set-CellAlign $ws.Range("B18") $xlCenter $xlGeneralV $false
Set-CellAlign $ws.Range("C18") $xlCenter $xlGeneralV $false
Set-CellAlign $ws.Range("D18") $xlCenter $xlGeneralV $false

Set-CellAlign $ws.Range("A19") $xlCenter $xlVCenter $true
Set-CellAlign $ws.Range("D19") $xlCenter $xlVCenter $false

Set-CellAlign $ws.Range("A20") $xlCenter $xlVCenter $true
Set-CellAlign $ws.Range("C20") $xlLeft $xlTop $false
Set-CellAlign $ws.Range("D20") $xlCenter $xlVCenter $false

# ---------------------------------------------------------------------
# 4. Row heights for the two data rows of the first test case.
# ---------------------------------------------------------------------
$ws.Rows.Item(14).RowHeight = 16.5
$ws.Rows.Item(15).RowHeight = 45

# Rows that used to carry a custom height (old row 6) go back to the
# sheet default now that their content/merge has moved.
$ws.Rows.Item(6).AutoFit()

# ---------------------------------------------------------------------
# 5. Merged cells for the new layout.
# ---------------------------------------------------------------------
$ws.Range("C14:C15").Merge()
$ws.Range("D14:D15").Merge()
$ws.Range("D19:D20").Merge()

# ---------------------------------------------------------------------
# 6. Selection, matching the saved cursor position in the target file.
# ---------------------------------------------------------------------
$ws.Range("E13").Select()
